# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (commit: "Updated cryptos list ... with GitHub Actions").
#
# Column D ("Price") values are free-form strings scraped from the site
# (e.g. "30.524.73", "1.000", "0.4798"). Many of them look like plain
# numbers to Excel, so before writing them we force the cell's number
# format to Text ("@") to guarantee they remain stored exactly as typed
# (preserving trailing zeros / thousands-grouped dots) instead of being
# auto-converted into a numeric value.
#
# Column E ("Volume(1h)") values are percentage strings that already
# contain non-numeric characters (%, spaces) so Excel keeps them as text
# automatically; no special number formatting is required for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-18: price/volume updates ---
$ws.Range("D2").Value = "30.524.73"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.918.82"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.47"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4798"
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2894"
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.22"
$ws.Range("E10").Value = "  +3.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.13"
$ws.Range("E11").Value = "  +4.78%  "
$ws.Range("D12").Value = "1.913.54"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07567"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.264"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6683"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "300.06"
$ws.Range("E16").Value = "  +3.00%  "
$ws.Range("D17").Value = "30.534.53"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.06"
$ws.Range("E18").Value = "  +0.98%  "

# --- Rows 19-20: Uniswap/Dai rows swap places (content-wise) ---
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.576"
$ws.Range("E20").Value = "  +5.84%  "

# --- Rows 21-51: price/volume updates ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007574"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "2.169.50"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.439"
$ws.Range("E24").Value = "  +3.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.488"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.48"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.33"
$ws.Range("E27").Value = "  -4.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.118"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("E29").Value = "  +0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.394"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.165"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.039"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04994"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7361"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9997"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02033"
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.684"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.85"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.018"
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4430"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.95"
$ws.Range("E43").Value = "  +6.69%  "
$ws.Range("E44").Value = "  -1.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.914"
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.0000"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.50"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.265"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.305"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1232"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2537"
$ws.Range("E51").Value = "  +3.23%  "
